$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "D" (Price) values are plain decimal numbers (e.g. "231.62").
# Assigning those directly to .Value would let Excel auto-convert them to
# numbers (losing formatting like trailing zeros, e.g. "1.90" -> 1.9).
# Force those specific cells to Text format first so the literal string is kept.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Updated prices (column D) and 1h volume-change percentages (column E)
$ws.Range("D2").Value = "34.912.66"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.840.17"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "231.62"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "39.87"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "0.328"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").Value = "0.0687"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "0.0986"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "2.108.18"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "11.44"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").Value = "1.844.22"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "0.673"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "34.926.59"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "69.89"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "240.53"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").Value = "4.69"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "2.28"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").Value = "171.38"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "7.79"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("D29").Value = "1.52"
$ws.Range("E29").Value = "  -5.31%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "3.95"
$ws.Range("E32").Value = "  -5.29%  "
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  +7.31%  "
$ws.Range("D35").Value = "1.24"
$ws.Range("E35").Value = "  +8.31%  "
$ws.Range("E36").Value = "  +11.64%  "
$ws.Range("D37").Value = "0.696"
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("E38").Value = "  +6.73%  "
$ws.Range("D39").Value = "91.10"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "1.342.62"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "14.81"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("D43").Value = "2.28"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").Value = "2.023.04"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "3.38"
$ws.Range("E50").Value = "  +18.56%  "
$ws.Range("E51").Value = "  +1.84%  "
